$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 95;  A = 94; B = "Monday, Jan 09"; C = "2:15 PM";  D = "LH1357"; E = "Frankfurt";  F = "(FRA)"; G = "Lufthansa ";                        H = "CRJ9"; I = "(D-ACNJ)"; J = "2:43 PM"; L = "0 hours, 28 minutes" },
    @{ Row = 96;  A = 95; B = "Monday, Jan 09"; C = "3:10 PM";  D = "LO3884"; E = "Warsaw";      F = "(WAW)"; G = "LOT (Sliwka Naleczowska Livery) "; H = "E195"; I = "(SP-LNC)";  J = "3:52 PM"; L = "0 hours, 42 minutes" },
    @{ Row = 97;  A = 96; B = "Monday, Jan 09"; C = "3:25 PM";  D = "W61283"; E = "Tel Aviv";    F = "(TLV)"; G = "Wizz Air ";                        H = "A321"; I = "(HA-LXE)"; J = "3:35 PM"; L = "0 hours, 10 minutes" },
    @{ Row = 98;  A = 97; B = "Monday, Jan 09"; C = "3:30 PM";  D = "LO6317"; E = "Punta Cana";  F = "(PUJ)"; G = "LOT ";                              H = "B788"; I = "(SP-LRC)"; J = "3:41 PM"; L = "0 hours, 11 minutes" },
    @{ Row = 99;  A = 98; B = "Monday, Jan 09"; C = "4:20 PM";  D = "W61167"; E = "Stavanger";   F = "(SVG)"; G = "Wizz Air ";                        H = "A321"; I = "(HA-LTC)"; J = "4:31 PM"; L = "0 hours, 11 minutes" },
    @{ Row = 100; A = 99; B = "Monday, Jan 09"; C = "4:50 PM";  D = "KL1816"; E = "Amsterdam";   F = "(AMS)"; G = "KLM ";                              H = "E295"; I = "(PH-NXF)"; J = "4:56 PM"; L = "0 hours, 6 minutes" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 12).Value = $r.L
}
